$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "09:00am-05:30pm"
$ws.Range("B9").Value = "09:00am-05:30pm"
$ws.Range("B11").Value = "09:00am-05:30pm"
$ws.Range("B13").Value = "09:00am-05:30pm"
$ws.Range("B15").Value = "09:00am-05:30pm"
